# Fixed puzzle issues raised during dry run
# (mala_stalls menu: correct the "cold" price column, H10:H14)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu")

$ws.Range("H10").Value = 0
$ws.Range("H12").Value = 0.4
$ws.Range("H13").Value = 0.6
$ws.Range("H14").Value = 1

# Drop the stray selection that had been left on H16 by resetting it to A1
$ws.Range("A1").Select()
